# Update the cryptos list (price + 1h volume%) with latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume 1h) hold text values (e.g. "66.916.40",
# "  +0.57%  "). Force the range to Text format first so Excel does not
# reinterpret numeric-looking strings (like "594.10") as numbers and drop
# the trailing zero / thousands-style dots.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "66.916.40"
$ws.Range("E2").Value = "  +0.57%  "

$ws.Range("D3").Value = "3.497.27"
$ws.Range("E3").Value = "  +0.33%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "594.10"
$ws.Range("E5").Value = "  +0.44%  "

$ws.Range("D6").Value = "172.64"
$ws.Range("E6").Value = "  +2.17%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("D9").Value = "0.131"
$ws.Range("E9").Value = "  +3.41%  "

$ws.Range("D10").Value = "7.19"
$ws.Range("E10").Value = "  -1.86%  "

$ws.Range("E11").Value = "  -0.87%  "

$ws.Range("D12").Value = "4.102.45"
$ws.Range("E12").Value = "  +0.37%  "

$ws.Range("E13").Value = "  -0.24%  "

$ws.Range("D14").Value = "29.21"
$ws.Range("E14").Value = "  +3.91%  "

$ws.Range("D15").Value = "66.924.35"
$ws.Range("E15").Value = "  +0.51%  "

$ws.Range("E16").Value = "  +0.61%  "

$ws.Range("D17").Value = "3.488.48"
$ws.Range("E17").Value = "  +0.13%  "

$ws.Range("D18").Value = "6.28"
$ws.Range("E18").Value = "  -0.44%  "

$ws.Range("D19").Value = "14.26"
$ws.Range("E19").Value = "  +1.82%  "

$ws.Range("D20").Value = "394.53"
$ws.Range("E20").Value = "  +0.70%  "

$ws.Range("D21").Value = "7.95"
$ws.Range("E21").Value = "  +0.54%  "

$ws.Range("D22").Value = "73.47"

$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("D24").Value = "0.535"
$ws.Range("E24").Value = "  +0.23%  "

$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("E26").Value = "  -0.77%  "

$ws.Range("E27").Value = "  +0.44%  "

$ws.Range("E28").Value = "  -0.29%  "

$ws.Range("D29").Value = "6.18"
$ws.Range("E29").Value = "  -2.10%  "

$ws.Range("D30").Value = "1.42"
$ws.Range("E30").Value = "  -2.30%  "

$ws.Range("E31").Value = "  -0.26%  "

$ws.Range("D32").Value = "23.69"
$ws.Range("E32").Value = "  +0.43%  "

$ws.Range("D33").Value = "7.36"
$ws.Range("E33").Value = "  -0.40%  "

$ws.Range("D34").Value = "1.61"
$ws.Range("E34").Value = "  +0.42%  "

$ws.Range("D35").Value = "162.69"
$ws.Range("E35").Value = "  +0.67%  "

$ws.Range("E36").Value = "  -1.24%  "

$ws.Range("E37").Value = "  -0.35%  "

$ws.Range("E38").Value = "  +1.80%  "

$ws.Range("D39").Value = "4.64"
$ws.Range("E39").Value = "  +0.35%  "

$ws.Range("D40").Value = "0.0739"
$ws.Range("E40").Value = "  -0.75%  "

$ws.Range("D41").Value = "2.835.24"
$ws.Range("E41").Value = "  +2.59%  "

$ws.Range("D42").Value = "27.03"
$ws.Range("E42").Value = "  +0.71%  "

$ws.Range("D43").Value = "26.16"
$ws.Range("E43").Value = "  -1.21%  "

$ws.Range("D44").Value = "42.75"
$ws.Range("E44").Value = "  -0.89%  "

$ws.Range("E45").Value = "  +2.44%  "

$ws.Range("E46").Value = "  -3.01%  "

$ws.Range("D47").Value = "336.86"
$ws.Range("E47").Value = "  -2.58%  "

$ws.Range("D48").Value = "34.53"
$ws.Range("E48").Value = "  +1.84%  "

$ws.Range("D49").Value = "1.07"
$ws.Range("E49").Value = "  -1.27%  "

$ws.Range("E50").Value = "  -1.31%  "

$ws.Range("D51").Value = "0.841"
$ws.Range("E51").Value = "  -4.35%  "
